# Rename the sheet: "Sheet1" -> "Statistics"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Statistics"

# The two header cells (A1:B1) lose their bold/bordered/centered style —
# reset them back to the workbook's default "Normal" formatting.
$ws.Range("A1:B1").ClearFormats()

# B5 used to hold a (broken, locale-mismatched) SUM formula; replace it
# with the (also locale-mismatched) AVERAGE formula the author switched to.
$ws.Range("B5").Formula = "=PROMEDIO(B2:B4)"

# A5 was an empty placeholder cell (inline string with no text) — drop it
# entirely so the row only contains the B5 formula cell.
$ws.Range("A5").ClearContents()

# Tighten the page margins to the (smaller) Excel "Normal" defaults.
$ps = $ws.PageSetup
$ps.LeftMargin = 50.4
$ps.RightMargin = 50.4
$ps.TopMargin = 54
$ps.BottomMargin = 54
$ps.HeaderMargin = 21.599999999999998
$ps.FooterMargin = 21.599999999999998
